# Add new job-book rows 55-57 to the JOBS sheet.
# Columns: A jobNumber(number) B modelNumber C serialNumber D voltage
#          E unloaders F statorStatus G incomingNumber H scrap I notes
#          J enteredBy K enteredOn L _isDeleted(bool) M deletedBy
#          N deletedOn O warranty
# All columns except A (number) and L (boolean) are text in this sheet,
# even when the text looks like a number or a date, so any cell whose
# literal value could be auto-interpreted as a number/date by Excel is
# pre-formatted as Text ("@") before the value is written, forcing it to
# stay a text value instead of being silently converted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $text
}

function Set-PlainText($addr, $text) {
    $ws.Range($addr).Value = $text
}

# --- Row 55 ---
$ws.Range("A55").Value = 71307
Set-PlainText "B55" "4MK1 35X AWM D P"
Set-PlainText "C55" "21D 57594 M"
Set-PlainText "D55" "380 420YY"
Set-TextCell  "E55" "0"
Set-PlainText "F55" "?"
Set-TextCell  "G55" "82915"
Set-PlainText "H55" "NO"
Set-PlainText "I55" ""
Set-PlainText "J55" "ravi"
Set-TextCell  "K55" "9/8/2022"
$ws.Range("L55").Value = $false
Set-PlainText "M55" "N/A"
Set-PlainText "N55" "N/A"
Set-PlainText "O55" "NO"

# --- Row 56 ---
$ws.Range("A56").Value = 71308
Set-PlainText "B56" "O6DE5379DC1900"
Set-PlainText "C56" "3609UD8512"
Set-PlainText "D56" "MULTI"
Set-TextCell  "E56" "2"
Set-PlainText "F56" "GOOD"
Set-TextCell  "G56" "82939"
Set-PlainText "H56" "NO"
Set-PlainText "I56" ""
Set-PlainText "J56" "ravi"
Set-TextCell  "K56" "9/8/2022"
$ws.Range("L56").Value = $false
Set-PlainText "M56" "N/A"
Set-PlainText "N56" "N/A"
Set-PlainText "O56" "NO"

# --- Row 57 ---
$ws.Range("A57").Value = 71309
Set-PlainText "B57" "O6DG5373DC0600"
Set-PlainText "C57" "4307U00575"
Set-TextCell  "D57" "460"
Set-TextCell  "E57" "2"
Set-PlainText "F57" "?"
Set-TextCell  "G57" "82940"
Set-PlainText "H57" "NO"
Set-PlainText "I57" "2 s HD e unl"
Set-PlainText "J57" "ravi"
Set-TextCell  "K57" "9/8/2022"
$ws.Range("L57").Value = $false
Set-PlainText "M57" "N/A"
Set-PlainText "N57" "N/A"
Set-PlainText "O57" "NO"
